$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing "Coach" rater entries to "Coach1" ---
$ws.Range("D2").Value = "Coach1"
$ws.Range("D3").Value = "Coach1"

# --- Row 5: new "Coach2" evaluation for Player1 / End of Fall 2018 ---
$ws.Range("A5").Value = "Player1"
$ws.Range("B5").Value = 2003
$ws.Range("C5").Value = "End of Fall 2018"
$ws.Range("D5").Value = "Coach2"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 3
$ws.Range("G5").Value = 4
$ws.Range("H5").Value = 5
$ws.Range("I5").Value = 2
$ws.Range("J5").Value = 3
$ws.Range("K5").Value = 4
$ws.Range("L5").Value = 5
$ws.Range("M5").Value = 2
$ws.Range("N5").Value = 3
$ws.Range("O5").Value = 4
$ws.Range("P5").Value = 5
$ws.Range("Q5").Value = 2
$ws.Range("R5").Value = 2
$ws.Range("S5").Value = 3
$ws.Range("T5").Value = 3
$ws.Range("U5").Value = 4
$ws.Range("V5").Value = 4
$ws.Range("W5").Value = 5
$ws.Range("X5").Value = 5
$ws.Range("Y5").Value = 4
$ws.Range("Z5").Value = 4
$ws.Range("AA5").Value = 3
$ws.Range("AB5").Value = 3
$ws.Range("AC5").Value = 2
$ws.Range("AD5").Value = 2
$ws.Range("AE5").Value = 1
$ws.Range("AF5").Value = 1
$ws.Range("AG5").Value = 2
$ws.Range("AH5").Value = 3
$ws.Range("AI5").Value = 4
$ws.Range("AJ5").Value = 5
$ws.Range("AK5").Value = "I think he dislikes hotdogs"

# --- Row 6: new "Coach2" evaluation for Player2 / End of Fall 2019 ---
$ws.Range("A6").Value = "Player2"
$ws.Range("B6").Value = 2004
$ws.Range("C6").Value = "End of Fall 2019"
$ws.Range("D6").Value = "Coach2"
$ws.Range("E6").Value = 4
$ws.Range("F6").Value = 4
$ws.Range("G6").Value = 3
$ws.Range("H6").Value = 4
$ws.Range("I6").Value = 4
$ws.Range("J6").Value = 3
$ws.Range("K6").Value = 4
$ws.Range("L6").Value = 4
$ws.Range("M6").Value = 3
$ws.Range("N6").Value = 5
$ws.Range("O6").Value = 4
$ws.Range("P6").Value = 4
$ws.Range("Q6").Value = 3
$ws.Range("R6").Value = 2
$ws.Range("S6").Value = 3
$ws.Range("T6").Value = 4
$ws.Range("U6").Value = 2
$ws.Range("V6").Value = 3
$ws.Range("W6").Value = 4
$ws.Range("X6").Value = 2
$ws.Range("Y6").Value = 1
$ws.Range("Z6").Value = 3
$ws.Range("AA6").Value = 2
$ws.Range("AB6").Value = 4
$ws.Range("AC6").Value = 3
$ws.Range("AD6").Value = 5
$ws.Range("AE6").Value = 5
$ws.Range("AF6").Value = 3
$ws.Range("AG6").Value = 2
$ws.Range("AH6").Value = 1
$ws.Range("AI6").Value = 3
$ws.Range("AJ6").Value = 3
$ws.Range("AK6").Value = "Giggles when he dribbles"

# --- Restore view/selection state to match the post-edit selection ---
$ws.Activate() | Out-Null
$ws.Range("A5:XFD6").Select() | Out-Null
$ws.Range("AJ6").Select() | Out-Null
